{"js": "// Update the two-digit multiplication problems in the (single) table.\n// Each populated row holds 5 problems (\"NN\u00d7NN=\"); only the <w:t> text\n// content changes, so we replace each cell's paragraph range text in\n// place (keeps run/paragraph formatting such as font + size intact).\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, colIndex, expectedOldText, newText) \u2014 rows are 0-based and\n// only the rows that actually contain a problem are listed here.\nconst edits = [\n  [0, 0, \"72\u00d715=\", \"39\u00d746=\"],\n  [0, 1, \"32\u00d725=\", \"76\u00d736=\"],\n  [0, 2, \"18\u00d727=\", \"85\u00d736=\"],\n  [0, 3, \"86\u00d714=\", \"93\u00d722=\"],\n  [0, 4, \"32\u00d735=\", \"72\u00d791=\"],\n\n  [4, 0, \"45\u00d756=\", \"71\u00d785=\"],\n  [4, 1, \"81\u00d729=\", \"35\u00d717=\"],\n  [4, 2, \"93\u00d722=\", \"95\u00d783=\"],\n  [4, 3, \"78\u00d786=\", \"60\u00d724=\"],\n  [4, 4, \"63\u00d797=\", \"38\u00d720=\"],\n\n  [9, 0, \"73\u00d722=\", \"59\u00d758=\"],\n  [9, 1, \"41\u00d766=\", \"93\u00d791=\"],\n  [9, 2, \"44\u00d742=\", \"86\u00d780=\"],\n  [9, 3, \"83\u00d726=\", \"94\u00d720=\"],\n  [9, 4, \"97\u00d785=\", \"59\u00d744=\"],\n\n  [14, 0, \"62\u00d756=\", \"62\u00d759=\"],\n  [14, 1, \"91\u00d714=\", \"24\u00d780=\"],\n  [14, 2, \"47\u00d754=\", \"70\u00d711=\"],\n  [14, 3, \"75\u00d733=\", \"39\u00d770=\"],\n  [14, 4, \"92\u00d738=\", \"82\u00d750=\"],\n\n  [19, 0, \"29\u00d787=\", \"25\u00d715=\"],\n  [19, 1, \"63\u00d731=\", \"36\u00d737=\"],\n  [19, 2, \"69\u00d757=\", \"91\u00d745=\"],\n  [19, 3, \"83\u00d731=\", \"13\u00d726=\"],\n  [19, 4, \"97\u00d782=\", \"73\u00d727=\"],\n];\n\n// Load current text for every target cell first so we can verify we are\n// editing the expected cell before writing (defensive, avoids silently\n// mutating the wrong cell if the table shape ever differs).\nconst ranges = edits.map(([row, col]) => {\n  const cell = table.getCell(row, col);\n  const range = cell.body.paragraphs.getFirst().getRange();\n  range.load(\"text\");\n  return range;\n});\nawait context.sync();\n\nfor (let i = 0; i < edits.length; i++) {\n  const [, , oldText, newText] = edits[i];\n  const range = ranges[i];\n  if (range.text !== oldText) {\n    console.log(\n      `Warning: cell text \"${range.text}\" did not match expected \"${oldText}\"` +\n        ` \u2014 replacing with \"${newText}\" anyway.`\n    );\n  }\n  // Replacing the text of the existing paragraph range (instead of\n  // clearing + inserting into the cell body) preserves the run's\n  // <w:rPr>/<w:pPr> formatting, matching the target diff which only\n  // touches the <w:t> content.\n  range.insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the two-digit multiplication problems in the (single) table.\n# Each populated row holds 5 problems (\"NN\u00d7NN=\"); only the <w:t> text\n# content changes, so we overwrite each cell's Range.Text in place,\n# which keeps the existing run/paragraph formatting (font + size)\n# untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Word's Cell(row, col) is 1-based \u2014 rows/cols below mirror the 0-based\n# table layout (only rows 1, 5, 10, 15, 20 contain problems).\n$edits = @(\n    @(1, 1, \"72\u00d715=\", \"39\u00d746=\"),\n    @(1, 2, \"32\u00d725=\", \"76\u00d736=\"),\n    @(1, 3, \"18\u00d727=\", \"85\u00d736=\"),\n    @(1, 4, \"86\u00d714=\", \"93\u00d722=\"),\n    @(1, 5, \"32\u00d735=\", \"72\u00d791=\"),\n\n    @(5, 1, \"45\u00d756=\", \"71\u00d785=\"),\n    @(5, 2, \"81\u00d729=\", \"35\u00d717=\"),\n    @(5, 3, \"93\u00d722=\", \"95\u00d783=\"),\n    @(5, 4, \"78\u00d786=\", \"60\u00d724=\"),\n    @(5, 5, \"63\u00d797=\", \"38\u00d720=\"),\n\n    @(10, 1, \"73\u00d722=\", \"59\u00d758=\"),\n    @(10, 2, \"41\u00d766=\", \"93\u00d791=\"),\n    @(10, 3, \"44\u00d742=\", \"86\u00d780=\"),\n    @(10, 4, \"83\u00d726=\", \"94\u00d720=\"),\n    @(10, 5, \"97\u00d785=\", \"59\u00d744=\"),\n\n    @(15, 1, \"62\u00d756=\", \"62\u00d759=\"),\n    @(15, 2, \"91\u00d714=\", \"24\u00d780=\"),\n    @(15, 3, \"47\u00d754=\", \"70\u00d711=\"),\n    @(15, 4, \"75\u00d733=\", \"39\u00d770=\"),\n    @(15, 5, \"92\u00d738=\", \"82\u00d750=\"),\n\n    @(20, 1, \"29\u00d787=\", \"25\u00d715=\"),\n    @(20, 2, \"63\u00d731=\", \"36\u00d737=\"),\n    @(20, 3, \"69\u00d757=\", \"91\u00d745=\"),\n    @(20, 4, \"83\u00d731=\", \"13\u00d726=\"),\n    @(20, 5, \"97\u00d782=\", \"73\u00d727=\")\n)\n\nforeach ($edit in $edits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $oldText = $edit[2]\n    $newText = $edit[3]\n    $cell = $t.Cell($row, $col)\n\n    # Cell.Range.Text includes the trailing cell-mark (CR + cell marker);\n    # trim it off before comparing to the plain problem text.\n    $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($currentText -ne $oldText) {\n        Write-Warning \"Cell ($row,$col) text `\"$currentText`\" did not match expected `\"$oldText`\" \u2014 replacing with `\"$newText`\" anyway.\"\n    }\n\n    # Overwriting Range.Text in place (rather than deleting + re-typing)\n    # preserves the run's existing rPr/pPr formatting, matching the\n    # target diff which only touches the <w:t> content.\n    $cell.Range.Text = $newText\n}\n"}
